# Rebuild the Il7-Il7r edge table with updated TPM-derived values.
# New data covers Sending cluster x Target cluster in
# {ECs, FAPs, Inflammatory-Mac, MuSCs} x {Inflammatory-Mac, MuSCs, Neutrophils, Resolving-Mac}
# i.e. 16 data rows (row 2 .. row 17), replacing the previous 15-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,20

$row = @("ECs", "Il7", "Il7r", "Inflammatory-Mac", 1, 0.5, 0.3296565, 0.659313, 0.2343734281062221, 0.1802644281836906, 3, 1, 29.63884466666667, 88.916534, 0.2980837622946661, 0.2982355126435957, 9.770637796857001, 58.623826781142, 0.0698629132318011, 0.0537612541507676)
for ($c = 0; $c -lt 20; $c++) { $data[0, $c] = $row[$c] }
$row = @("ECs", "Il7", "Il7r", "MuSCs", 1, 0.5, 0.3296565, 0.659313, 0.2343734281062221, 0.1802644281836906, 1, 0.5, 0.15178, 0.30356, 0.001526481681385076, 0.00101817252816096, 0.05003526357, 0.20014105428, 0.00035776674460757, 0.000183540288581278)
for ($c = 0; $c -lt 20; $c++) { $data[1, $c] = $row[$c] }
$row = @("ECs", "Il7", "Il7r", "Neutrophils", 1, 0.5, 0.3296565, 0.659313, 0.2343734281062221, 0.1802644281836906, 3, 1, 9.907502666666666, 29.722508, 0.09964172703215479, 0.09969245326671611, 3.266072652834, 19.596435917004, 0.02335337314695054, 0.01797100308235387)
for ($c = 0; $c -lt 20; $c++) { $data[2, $c] = $row[$c] }
$row = @("ECs", "Il7", "Il7r", "Resolving-Mac", 1, 0.5, 0.3296565, 0.659313, 0.2343734281062221, 0.1802644281836906, 3, 1, 59.73313466666667, 179.199404, 0.600748028991794, 0.6010538615615272, 19.691416108242, 118.148496649452, 0.1407993749828629, 0.1083486306619878)
for ($c = 0; $c -lt 20; $c++) { $data[3, $c] = $row[$c] }
$row = @("FAPs", "Il7", "Il7r", "Inflammatory-Mac", 3, 1, 0.7595033333333333, 2.27851, 0.5399784317659018, 0.6229731588195906, 3, 1, 29.63884466666667, 88.916534, 0.2980837622946661, 0.2982355126435957, 22.51080132048222, 202.59721188434, 0.1609588024987536, 0.1857927193837608)
for ($c = 0; $c -lt 20; $c++) { $data[4, $c] = $row[$c] }
$row = @("FAPs", "Il7", "Il7r", "MuSCs", 3, 1, 0.7595033333333333, 2.27851, 0.5399784317659018, 0.6229731588195906, 1, 0.5, 0.15178, 0.30356, 0.001526481681385076, 0.00101817252816096, 0.1152774159333333, 0.6916644956, 0.0008242671844336902, 0.0006342941560917616)
for ($c = 0; $c -lt 20; $c++) { $data[5, $c] = $row[$c] }
$row = @("FAPs", "Il7", "Il7r", "Neutrophils", 3, 1, 0.7595033333333333, 2.27851, 0.5399784317659018, 0.6229731588195906, 3, 1, 9.907502666666666, 29.722508, 0.09964172703215479, 0.09969245326671611, 7.524781300342221, 67.72303170307998, 0.053804383501269, 0.06210572252204054)
for ($c = 0; $c -lt 20; $c++) { $data[6, $c] = $row[$c] }
$row = @("FAPs", "Il7", "Il7r", "Resolving-Mac", 3, 1, 0.7595033333333333, 2.27851, 0.5399784317659018, 0.6229731588195906, 3, 1, 59.73313466666667, 179.199404, 0.600748028991794, 0.6010538615615272, 45.36751488978222, 408.30763400804, 0.3243909785814454, 0.3744404227576975)
for ($c = 0; $c -lt 20; $c++) { $data[7, $c] = $row[$c] }
$row = @("Inflammatory-Mac", "Il7", "Il7r", "Inflammatory-Mac", 1, 0.3333333333333333, 0.084886, 0.254658, 0.06035076759664913, 0.06962668528059096, 3, 1, 29.63884466666667, 88.916534, 0.2980837622946661, 0.2982355126435957, 2.515922968374667, 22.643306715372, 0.0179895838625802, 0.02076515017833135)
for ($c = 0; $c -lt 20; $c++) { $data[8, $c] = $row[$c] }
$row = @("Inflammatory-Mac", "Il7", "Il7r", "MuSCs", 1, 0.3333333333333333, 0.084886, 0.254658, 0.06035076759664913, 0.06962668528059096, 1, 0.5, 0.15178, 0.30356, 0.001526481681385076, 0.00101817252816096, 0.01288399708, 0.07730398247999999, [double]"9.212434119381292E-05", [double]"7.089197817960677E-05")
for ($c = 0; $c -lt 20; $c++) { $data[9, $c] = $row[$c] }
$row = @("Inflammatory-Mac", "Il7", "Il7r", "Neutrophils", 1, 0.3333333333333333, 0.084886, 0.254658, 0.06035076759664913, 0.06962668528059096, 3, 1, 9.907502666666666, 29.722508, 0.09964172703215479, 0.09969245326671611, 0.8410082713626666, 7.569074442263999, 0.006013454711046326, 0.006941255068451665)
for ($c = 0; $c -lt 20; $c++) { $data[10, $c] = $row[$c] }
$row = @("Inflammatory-Mac", "Il7", "Il7r", "Resolving-Mac", 1, 0.3333333333333333, 0.084886, 0.254658, 0.06035076759664913, 0.06962668528059096, 3, 1, 59.73313466666667, 179.199404, 0.600748028991794, 0.6010538615615272, 5.070506869314667, 45.634561823832, 0.0362556046818288, 0.04184938805562834)
for ($c = 0; $c -lt 20; $c++) { $data[11, $c] = $row[$c] }
$row = @("MuSCs", "Il7", "Il7r", "Inflammatory-Mac", 1, 0.5, 0.232498, 0.464996, 0.1652973725312269, 0.1271357277161278, 3, 1, 29.63884466666667, 88.916534, 0.2980837622946661, 0.2982355126435957, 6.890972107310668, 41.345832643864, 0.04927246270153112, 0.03791638893073598)
for ($c = 0; $c -lt 20; $c++) { $data[12, $c] = $row[$c] }
$row = @("MuSCs", "Il7", "Il7r", "MuSCs", 1, 0.5, 0.232498, 0.464996, 0.1652973725312269, 0.1271357277161278, 1, 0.5, 0.15178, 0.30356, 0.001526481681385076, 0.00101817252816096, 0.03528854644, 0.14115418576, 0.0002523234111500025, 0.0001294461053083133)
for ($c = 0; $c -lt 20; $c++) { $data[13, $c] = $row[$c] }
$row = @("MuSCs", "Il7", "Il7r", "Neutrophils", 1, 0.5, 0.232498, 0.464996, 0.1652973725312269, 0.1271357277161278, 3, 1, 9.907502666666666, 29.722508, 0.09964172703215479, 0.09969245326671611, 2.303474554994667, 13.820847329968, 0.01647051567288892, 0.01267447259387002)
for ($c = 0; $c -lt 20; $c++) { $data[14, $c] = $row[$c] }
$row = @("MuSCs", "Il7", "Il7r", "Resolving-Mac", 1, 0.5, 0.232498, 0.464996, 0.1652973725312269, 0.1271357277161278, 3, 1, 59.73313466666667, 179.199404, 0.600748028991794, 0.6010538615615272, 13.88783434373067, 83.32700606238402, 0.0993020707456569, 0.07641542008621351)
for ($c = 0; $c -lt 20; $c++) { $data[15, $c] = $row[$c] }

$ws.Range("A2:T17").Value = $data

# Dimension (A1:T16 -> A1:T17) and shared strings are updated automatically by Excel.
